$d = $word.ActiveDocument

# 1) Update the placeholder ID text in the first paragraph.
#    Replaces the two runs ("**ID__AFFARS_5328_topic_5__ID**" + trailing
#    space) with a single run containing the new ID text.
$d.Content.Find.Execute(
    "**ID__AFFARS_5328_topic_5__ID** ", $true, $false, $false, $false,
    $false, $true, 1, $false, "**ID__AFFARS_5328_106_2__ID**", 2) | Out-Null

# 2) Update the first paragraph's formatting: add paragraph border spacing
#    (top/left/bottom/right, 5 each) and change the left indent from
#    120 twips (6pt) to 225 twips (11.25pt).
$p1 = $d.Paragraphs.Item(1)
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5
$p1.Format.LeftIndent = 11.25
